{"js": "// The tables in this document each have a \"Date\" column whose values were\n// recomputed by an Excel SUM-based formula and shifted 3 days later\n// (28.06.2020 -> 01.07.2020, 29.06.2020 -> 02.07.2020, 30.06.2020 -> 03.07.2020,\n// 01.07.2020 -> 04.07.2020). Every occurrence of each old date string in the\n// document maps to the same new date, so we can do a global text replace per\n// date value.\n//\n// Because \"01.07.2020\" is simultaneously an OLD value (-> \"04.07.2020\") and\n// the NEW value produced by another replacement (\"28.06.2020\" -> \"01.07.2020\"),\n// we must gather every match for every old date BEFORE performing any of the\n// text replacements; otherwise a later replacement could re-match text that\n// an earlier replacement just inserted.\n\nconst dateMap = {\n  \"28.06.2020\": \"01.07.2020\",\n  \"29.06.2020\": \"02.07.2020\",\n  \"30.06.2020\": \"03.07.2020\",\n  \"01.07.2020\": \"04.07.2020\",\n};\n\n// 1) Search for every old date string first and load the hit ranges.\nconst pending = [];\nfor (const oldDate of Object.keys(dateMap)) {\n  const results = context.document.body.search(oldDate, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  pending.push({ oldDate, results });\n}\nawait context.sync();\n\n// 2) Now that all hits (against the original text) are known, replace the\n// text in each range with the corresponding new date.\nfor (const { oldDate, results } of pending) {\n  const newDate = dateMap[oldDate];\n  for (const range of results.items) {\n    range.insertText(newDate, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Each of the document's 3 tables has a \"Date\" column (3rd column) whose\n# values were recalculated by an Excel SUM formula and landed 3 days later\n# than before (e.g. 28.06.2020 -> 01.07.2020, 01.07.2020 -> 04.07.2020, ...).\n#\n# Rather than hard-coding the old/new string pairs (which would risk a\n# chained-replacement bug since \"01.07.2020\" is both an old value and the\n# result of another replacement), read each date cell's current value,\n# compute \"+3 days\" from it, and write the result straight back into that\n# same cell's Range. Setting Range.Text in place (instead of Find/Replace)\n# also keeps every paragraph/run formatting attribute (rsid, etc.) untouched,\n# exactly like the reference edit.\n\n$d = $word.ActiveDocument\n\nforeach ($t in $d.Tables) {\n    $rowCount = $t.Rows.Count\n    # Row 1 is the header (\"Caption\", \"Group\", \"Date\", \"SUM1\", \"SUM2\");\n    # data rows start at row 2. Column 3 is \"Date\".\n    for ($r = 2; $r -le $rowCount; $r++) {\n        $cell = $t.Cell($r, 3)\n        $rng = $cell.Range\n        # Cell ranges include trailing cell-mark characters (CR + BEL);\n        # strip them before parsing the visible text.\n        $clean = $rng.Text.TrimEnd([char]13, [char]7)\n\n        if ($clean -match '^(\\d{2})\\.(\\d{2})\\.(\\d{4})$') {\n            $day = [int]$matches[1]\n            $month = [int]$matches[2]\n            $year = [int]$matches[3]\n\n            $oldDate = Get-Date -Year $year -Month $month -Day $day\n            $newDate = $oldDate.AddDays(3)\n\n            $rng.Text = $newDate.ToString(\"dd.MM.yyyy\")\n        }\n    }\n}\n"}
